$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 ("Hydrogen" / Non-metallic minerals) previously held a numeric 0;
# the corrected results clear it out entirely.
$ws.Range("D3").Value = ""

# Row 7 used to be the catch-all "Other" row with the Biogas total
# (277.5542358471358). The corrected data renames it to "Biogas" and
# zeroes out its value, moving the real "Other" total to a new row 8.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 0

# Insert the new "Other" row (row 8), copying row 7's label formatting
# (bold, bordered, centered) onto the new label cell before setting values.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 277.5542358471358
